# The data rows (3-48) got shuffled: for each row, the values in columns
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), O (Origen) and P (Precio $/Kg) were
# replaced by the values that used to live in another row of the same
# block, per the mapping below (row -> source row it now carries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    3 = 10
    4 = 14
    5 = 46
    6 = 21
    7 = 8
    8 = 18
    9 = 48
    10 = 41
    11 = 34
    12 = 26
    13 = 44
    14 = 45
    15 = 22
    16 = 6
    17 = 30
    18 = 24
    19 = 5
    20 = 4
    21 = 40
    22 = 16
    23 = 17
    24 = 25
    25 = 11
    26 = 32
    27 = 35
    28 = 33
    29 = 47
    30 = 13
    31 = 42
    32 = 27
    33 = 37
    34 = 7
    35 = 31
    36 = 28
    37 = 38
    38 = 20
    39 = 15
    40 = 3
    41 = 43
    42 = 9
    43 = 29
    44 = 39
    45 = 36
    46 = 12
    47 = 19
    48 = 23
}

$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot the original values for every affected row/column before any
# writes happen, because the mapping is a permutation: row 3 reads from
# row 10 while row 10 itself is also being overwritten, etc.
$snapshot = @{}
foreach ($row in $map.Keys) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowData
}

foreach ($row in $map.Keys) {
    $srcRow = $map[$row]
    $srcData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $srcData[$col]
    }
}
